# Update the "想去人数" (want-to-go count) figures for a few events.
# Sheet "展览" (sheet1): F2 289->290, F4 1110->1115, F5 579->580
# Sheet "全部类型" (sheet4): F2 289->290, F4 1110->1115, F6 579->580

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 290
$wsExhibition.Range("F4").Value = 1115
$wsExhibition.Range("F5").Value = 580

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 290
$wsAll.Range("F4").Value = 1115
$wsAll.Range("F6").Value = 580
